$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "automotive/bitcount/bitcnts 100000000"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0

$ws.Columns.Item(1).ColumnWidth = 33

$ws.Range("A16").Select()
